# Applies the MitsosBarton2006Ex313 "alpha_zero" regeneration edit:
# updates the expression / evaluation / bound cells on the
# Restricciones_del_lider, Restricciones_del_follower, Punto_modificado,
# Vector_bf and Vector_BF sheets to the new generated values.

$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param(
        $Worksheet,
        [string]$Cell,
        [string]$Value
    )
    $rng = $Worksheet.Range($Cell)
    # Force text storage so purely-numeric strings (e.g. "-3.3000000000000003")
    # are kept as text instead of being coerced into numeric cells, matching
    # the original workbook where every value is stored as a shared string.
    $rng.NumberFormat = "@"
    $rng.Value = $Value
    $rng.ClearFormats()
}

# NOTE: worksheet names are looked up by index rather than by name because
# "Vector_bf" and "Vector_BF" differ only by case and Worksheets.Item(name)
# resolution is case-insensitive (both would otherwise resolve to the same
# sheet).

# --- Restricciones_del_lider ---
$wsLider = $wb.Worksheets.Item(2)
Set-TextValue $wsLider "A2" "2.3000000000000003 - x"
Set-TextValue $wsLider "B2" "-3.3000000000000003"
Set-TextValue $wsLider "D2" "0.51"
Set-TextValue $wsLider "A3" "-2.3000000000000003 + x"
Set-TextValue $wsLider "B3" "1.3000000000000003"
Set-TextValue $wsLider "D3" "0.17"

# --- Restricciones_del_follower ---
$wsFollower = $wb.Worksheets.Item(3)
Set-TextValue $wsFollower "A2" "-4.449999999999999 + y"
Set-TextValue $wsFollower "B2" "3.4499999999999993"
Set-TextValue $wsFollower "D2" "0.82"
Set-TextValue $wsFollower "E2" "3.2"
Set-TextValue $wsFollower "F2" "9.8"
Set-TextValue $wsFollower "A3" "4.449999999999999 - y"
Set-TextValue $wsFollower "B3" "-5.449999999999999"
Set-TextValue $wsFollower "D3" "0.81"
Set-TextValue $wsFollower "E3" "7.199999999999999"
Set-TextValue $wsFollower "F3" "1.7000000000000002"

# --- Punto_modificado ---
$wsPunto = $wb.Worksheets.Item(4)
Set-TextValue $wsPunto "A2" "2.3000000000000003"
Set-TextValue $wsPunto "B2" "4.449999999999999"

# --- Vector_bf (sheet index 5) ---
$wsVecBf = $wb.Worksheets.Item(5)
Set-TextValue $wsVecBf "A2" "1.9220000000000041"

# --- Vector_BF (sheet index 6) ---
$wsVecBF = $wb.Worksheets.Item(6)
Set-TextValue $wsVecBF "A2" "-0.66"
Set-TextValue $wsVecBF "A3" "4.999999999999999"
